# Arbeitszeiten.xlsx update
# - fix the year typo in rows 89-94 (2018 -> 2019, i.e. +365 days)
# - add two new log entries (rows 95 & 96) for the translation backend work
# - update the sheet view (scrolled down, new selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Fix the mis-typed dates in rows 89-94 (were 2018-01-xx, should be 2019-01-xx) ---
$ws.Range("A89").Value = 43467
$ws.Range("A90").Value = 43468
$ws.Range("A91").Value = 43469
$ws.Range("A92").Value = 43475
$ws.Range("A93").Value = 43475
$ws.Range("A94").Value = 43481

# --- Append the new rows for the translation service work ---
$ws.Range("A95").Value = 43499
$ws.Range("B95").Value = "David"
$ws.Range("C95").Value = "Translation.csv erstellen + Translation in Db einspeisen"
$ws.Range("D95").Value = 10

$ws.Range("A96").Value = 43500
$ws.Range("B96").Value = "David"
$ws.Range("C96").Value = "Algorithmus für Translation + Englisch Übersetzung"
$ws.Range("D96").Value = 7

# --- Update the visible view: scrolled to row 77, new selection C97 ---
$ws.Application.ActiveWindow.ScrollRow = 77
$ws.Range("C97").Select()
